$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[49.95559502467271, 50.045957545858464]"
$ws.Range("T2").Value = "[49.95198635904036, 50.01632831257626]"
$ws.Range("L3").Value = "[49.99486815256596, 50.09525034968921]"
$ws.Range("T3").Value = "[49.958091116594275, 50.023849135655254]"
